$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Switch workbook back to automatic calculation ---
$excel.Calculation = -4105

# --- Bold header row (A1:C1) ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1:C1").Font.Bold = $true

# --- Append newly loaded railway od_pairs (rows 2556-2581) ---
# Write column-by-column (A, then B, then C) to mirror the original
# authoring order used when the data set was produced.
$colA = @(
  "78-21",
  "66-21",
  "83-21",
  "64-21",
  "67-21",
  "62-21",
  "49-45",
  "63-21",
  "51-45",
  "49-21",
  "56-21",
  "72-21",
  "59-21",
  "52-45",
  "29-21",
  "27-21",
  "32-21",
  "60-21",
  "28-21",
  "58-21",
  "25-21",
  "26-21",
  "57-21",
  "46-45",
  "23-21",
  "24-21"
)
$colB = @(
  "078-083-1002-065-055-056-1004-021",
  "066-067-014-015-017-021",
  "083-064-063-014-057-021",
  "064-063-014-057-021",
  "067-014-015-017-021",
  "No tiene",
  "049-051-1052-1037-046-1046-045",
  "063-014-057-021",
  "051-1052-1037-046-1046-045",
  "049-095-1059-060-1023-023-021",
  "056-058-020-019-021",
  "No tiene",
  "059-058-020-019-021",
  "052-1052-1037-046-1046-045",
  "029-1023-023-021",
  "027-007-026-025-024-1022-021",
  "032-027-007-026-025-024-1022-021",
  "060-1023-023-021",
  "028-029-1023-023-021",
  "058-020-019-021",
  "025-024-1022-021",
  "026-025-024-1022-021",
  "057-021",
  "046-1046-045",
  "023-021",
  "024-1022-021"
)
$colC = @(
  "angosta",
  "angosta",
  "ancha",
  "ancha",
  "angosta",
  "No encontrada",
  "ancha",
  "ancha",
  "ancha",
  "ancha",
  "ancha",
  "No encontrada",
  "ancha",
  "ancha",
  "ancha",
  "ancha",
  "ancha",
  "ancha",
  "ancha",
  "ancha",
  "ancha",
  "ancha",
  "ancha",
  "ancha",
  "ancha",
  "ancha"
)

$startRow = 2556
for ($i = 0; $i -lt $colA.Length; $i++) {
  $ws.Cells.Item($startRow + $i, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt $colB.Length; $i++) {
  $ws.Cells.Item($startRow + $i, 2).Value = $colB[$i]
}
for ($i = 0; $i -lt $colC.Length; $i++) {
  $ws.Cells.Item($startRow + $i, 3).Value = $colC[$i]
}

# --- Two trailing blank rows from the source data (rows 2582-2583) ---
$ws.Rows.Item(2582).Font.Bold = $false
$ws.Rows.Item(2583).Font.Bold = $false

# --- Reset the view to the top-left cell (A1), clearing the old scroll/selection ---
$excel.Goto($ws.Range("A1"))
